# Updates cryptos list cell values (Price column D, Volume(1h) column E)
# matching the scraped coinranking.com snapshot for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.111.28'
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").Value = '2.482.28'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").ClearContents()
$ws.Range("D5").Value = "'585.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("D6").ClearContents()
$ws.Range("D6").Value = "'171.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.76%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -0.19%  '

$ws.Range("D9").Value = '2.481.80'
$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("E10").Value = '  +1.07%  '

$ws.Range("E11").Value = '  +0.13%  '

$ws.Range("D12").ClearContents()
$ws.Range("D12").Value = "'4.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("D13").ClearContents()
$ws.Range("D13").Value = "'0.331"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.52%  '

$ws.Range("D15").ClearContents()
$ws.Range("D15").Value = "'25.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.84%  '

$ws.Range("D16").Value = '66.969.58'
$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("E17").Value = '  -1.80%  '

$ws.Range("D18").Value = '2.480.33'
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("E19").Value = '  -4.50%  '

$ws.Range("D20").ClearContents()
$ws.Range("D20").Value = "'7.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.36%  '

$ws.Range("E21").Value = '  -3.36%  '

$ws.Range("E22").Value = '  -1.06%  '

$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("D24").ClearContents()
$ws.Range("D24").Value = "'68.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.10%  '

$ws.Range("D25").ClearContents()
$ws.Range("D25").Value = "'4.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.66%  '

$ws.Range("E26").Value = '  -2.22%  '

$ws.Range("D27").ClearContents()
$ws.Range("D27").Value = "'9.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.87%  '

$ws.Range("E28").Value = '  -1.34%  '

$ws.Range("D29").Value = '2.608.40'
$ws.Range("E29").Value = '  -0.90%  '

$ws.Range("D30").Value = '0.0₃0899'
$ws.Range("E30").Value = '  -2.76%  '

$ws.Range("D31").ClearContents()
$ws.Range("D31").Value = "'509.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.16%  '

$ws.Range("E32").Value = '  -5.08%  '

$ws.Range("E33").Value = '  -3.12%  '

$ws.Range("E34").Value = '  -3.61%  '

$ws.Range("D35").ClearContents()
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").ClearContents()
$ws.Range("D36").Value = "'159.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.04%  '

$ws.Range("E37").Value = '  -6.91%  '

$ws.Range("E38").Value = '  +0.71%  '

$ws.Range("E39").Value = '  -3.67%  '

$ws.Range("E40").Value = '  -5.97%  '

$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("E42").Value = '  -3.46%  '

$ws.Range("E43").Value = '  -2.82%  '

$ws.Range("E44").Value = '  -1.49%  '

$ws.Range("D45").ClearContents()
$ws.Range("D45").Value = "'2.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.08%  '

$ws.Range("D46").ClearContents()
$ws.Range("D46").Value = "'38.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.90%  '

$ws.Range("D47").ClearContents()
$ws.Range("D47").Value = "'142.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.19%  '

$ws.Range("E48").Value = '  -4.14%  '

$ws.Range("E49").Value = '  -4.28%  '

$ws.Range("E50").Value = '  -6.12%  '

$ws.Range("E51").Value = '  -1.26%  '
